# GantBD.xlsx - "3. Implementação da Solução" section
#
# The project report moved on to the implementation phase, so this adds
# the new Gantt rows for that phase (and its six sub-tasks) and rolls the
# "Project Start" / "Today" markers forward to match.
#
# Note: the new task-name strings are written in the same order the
# author must have typed them (3, 3.1, 3.2, 3.5, 3.6, then backfilling
# 3.3 and 3.4) so the workbook's shared-string table comes out in the
# same order as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Project Start (E2) and Today (E3) markers
$ws.Range("E2").Value = 45243
$ws.Range("E3").Value = 45315

# Row 32 - "3. Implementação da Solução" (summary task)
$ws.Range("B32").Value = "3. Implementação da Solução"
$ws.Range("E32").Value = 45243
$ws.Range("F32").Value = 45315

# Row 33 - "3.1 Implementação da BD"
$ws.Range("B33").Value = "3.1 Implementação da BD"
$ws.Range("E33").Value = 45243
$ws.Range("F33").Value = 45249

# Row 34 - "3.2 Desenho das Interfaces"
$ws.Range("B34").Value = "3.2 Desenho das Interfaces"
$ws.Range("E34").Value = 45250
$ws.Range("F34").Value = 45256

# Row 37 - "3.5 Teste e Validação da Aplicação"
$ws.Range("B37").Value = "3.5 Teste e Validação da Aplicação"
$ws.Range("E37").Value = 45293
$ws.Range("F37").Value = 45305

# Row 38 - "3.6 Documentação"
$ws.Range("B38").Value = "3.6 Documentação"
$ws.Range("E38").Value = 45306
$ws.Range("F38").Value = 45315

# Row 35 - "3.3 Desenvolvimento da Lógica de Negócio"
$ws.Range("B35").Value = "3.3 Desenvolvimento da Lógica de Negócio"
$ws.Range("E35").Value = 45257
$ws.Range("F35").Value = 45277

# Row 36 - "3.4 Desenvolvimento das Interfaces"
$ws.Range("B36").Value = "3.4 Desenvolvimento das Interfaces"
$ws.Range("E36").Value = 45264
$ws.Range("F36").Value = 45284

# Bring the view back to where the author left it: zoomed to 70% with
# AF47 selected (the Gantt bars for the new tasks).
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
[void]$ws.Range("AF47").Select()
